$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pollen section additions (written in the same order the strings were
# originally authored, so the shared-string table order matches too)
$ws.Range("C62").Value = "LM"
$ws.Range("D62").Value = "SEM"
$ws.Range("E62").Value = "TEM"
$ws.Range("C61").Value = "(feel free to collapse these into one list if you prefer)"
$ws.Range("C63").Value = "grain size"
$ws.Range("D64").Value = "pore size"
$ws.Range("D65").Value = "pore depth"
$ws.Range("D66").Value = "annulus width"
$ws.Range("D67").Value = "annulus depth"
$ws.Range("C66").Value = "presence/absence of operculum"
$ws.Range("D63").Value = "exine ornamentation"
$ws.Range("C65").Value = "physical pollen description"
$ws.Range("C64").Value = "grain shape"
$ws.Range("E63").Value = "pollen wall ultrastructure"
$ws.Range("E64").Value = "depth of exine ornamentation"
$ws.Range("E65").Value = "presence/absence of microchannels"
$ws.Range("E66").Value = "size of microchannels"
$ws.Range("E67").Value = "orientation of microchannels"

$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

$ws.Application.ActiveWindow.ScrollRow = 50
$ws.Range("F67").Select()
